# Fruta / hortaliza, semanal
# Insert a new weekly record at row 33 (shifting the existing rows 33..129
# down to 34..130) on the active sheet of the "Papa" (Potato) subset sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 33..129 down one row, leaving a blank row 33 for the new record.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new data point.
$ws.Range("A33").Value = 11
$ws.Range("B33").Value = 'Vega Monumental Concepción'
$ws.Range("C33").Value = 'Bíobío'
$ws.Range("D33").Value = 44497
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = 100114001
$ws.Range("G33").Value = 'Papa'
$ws.Range("H33").Value = 'Patagonia'
$ws.Range("I33").Value = '1a (guarda)'
$ws.Range("J33").Value = 270
$ws.Range("K33").Value = 9000
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = 9556
$ws.Range("N33").Value = '$/saco 25 kilos'
$ws.Range("O33").Value = 'Región de Los Lagos'
$ws.Range("P33").Value = 382
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = 'Hortaliza'
